# Added periodic & upfront related scenarios
# Update the "repaymentstrategy" value on the ProductLoanInput sheet (row 17)
# from the stale "RBI (India)" placeholder to the correct scenario text.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("ProductLoanInput")

$ws1.Range("B17").Value = "Overdue/Due Fee/Int,Principal"

# Make ProductLoanInput the active sheet/tab, with B17 as the selected cell
# (mirrors the author re-selecting this sheet/cell after editing it).
$ws1.Activate()
[void]$ws1.Range("B17").Select()
